$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 289 (pushes ADB/SFO/FSD rows down by one)
$ws.Rows.Item(289).Insert()

# Populate the new row with the Fukuoka, Japan colo entry
$a = $ws.Cells.Item(289, 1)
$a.Value = "FUK"
$a.Font.Bold = $true
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4160
$a.Borders.LineStyle = 1

$ws.Cells.Item(289, 2).Value = "Fukuoka, Japan"
$ws.Cells.Item(289, 3).Value = 33.5902
$ws.Cells.Item(289, 4).Value = 130.4017
$ws.Cells.Item(289, 5).Value = "JP"
$ws.Cells.Item(289, 6).Value = "Asia Pacific"
$ws.Cells.Item(289, 7).Value = "Fukuoka"
